$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 193
$ws.Range("I2").Value = 520
$ws.Range("J2").Value = 2188
$ws.Range("K2").Value = 15
$ws.Range("L2").Value = 588
$ws.Range("M2").Value = 36
$ws.Range("N2").Value = 365
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 8
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 25
$ws.Range("S2").Value = 212
$ws.Range("T2").Value = 384
$ws.Range("U2").Value = 30
$ws.Range("V2").Value = 3354
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 3239
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 54
$ws.Range("AA2").Value = 24
